# Scheduled runner update: refresh market-price-derived columns
# (currentAveragePrice / currentAveragePriceNQ / currentAveragePriceHQ /
#  LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ) on the
# per-job Leve profit sheets. Values only -- no formulas, no formatting.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 4124.091
$ws.Range("J70").Value = 4596.1113
$ws.Range("L70").Value = 13788.3339
$ws.Range("N70").Value = -14328.3339

$ws.Range("H73").Value = 4124.091
$ws.Range("J73").Value = 4596.1113
$ws.Range("L73").Value = 13788.3339
$ws.Range("N73").Value = -15660.3339

$ws.Range("H128").Value = 41846
$ws.Range("J128").Value = 41846
$ws.Range("L128").Value = 41846
$ws.Range("N128").Value = -51806

$ws.Range("H129").Value = 845.37
$ws.Range("I129").Value = 355.6
$ws.Range("J129").Value = 899.7889
$ws.Range("K129").Value = 1066.8
$ws.Range("L129").Value = 2699.3667
$ws.Range("M129").Value = 3933.2
$ws.Range("N129").Value = -12699.3667

$ws.Range("H132").Value = 723968.3
$ws.Range("I132").Value = 13061.444
$ws.Range("J132").Value = 2003600.6
$ws.Range("K132").Value = 39184.33199999999
$ws.Range("L132").Value = 6010801.800000001
$ws.Range("M132").Value = -36654.33199999999
$ws.Range("N132").Value = -6015861.800000001

$ws.Range("H137").Value = 2219.5483
$ws.Range("I137").Value = 1558.4286
$ws.Range("J137").Value = 3607.9
$ws.Range("K137").Value = 4675.2858
$ws.Range("L137").Value = 10823.7
$ws.Range("M137").Value = -2125.2858
$ws.Range("N137").Value = -15923.7

$ws.Range("H138").Value = 4592.77
$ws.Range("I138").Value = 783.087
$ws.Range("K138").Value = 2349.261
$ws.Range("M138").Value = 2790.739

$ws.Range("H141").Value = 5653.2617
$ws.Range("I141").Value = 5891.846
$ws.Range("J141").Value = 2551.6667
$ws.Range("K141").Value = 17675.538
$ws.Range("L141").Value = 7655.000100000001
$ws.Range("M141").Value = -12495.538
$ws.Range("N141").Value = -18015.0001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4389.7964
$ws.Range("I32").Value = 4212.4346
$ws.Range("K32").Value = 4212.4346
$ws.Range("M32").Value = -3925.4346

$ws.Range("H61").Value = 1239.9791
$ws.Range("I61").Value = 1236.081
$ws.Range("J61").Value = 1253.091
$ws.Range("K61").Value = 1236.081
$ws.Range("L61").Value = 1253.091
$ws.Range("M61").Value = -1024.081
$ws.Range("N61").Value = -1677.091

$ws.Range("H97").Value = 897.6429000000001
$ws.Range("I97").Value = 560.1111
$ws.Range("J97").Value = 10011
$ws.Range("K97").Value = 560.1111
$ws.Range("L97").Value = 10011
$ws.Range("M97").Value = -64.11109999999996
$ws.Range("N97").Value = -11003

$ws.Range("H102").Value = 2740.5
$ws.Range("I102").Value = 2077.3333
$ws.Range("J102").Value = 3403.6667
$ws.Range("K102").Value = 2077.3333
$ws.Range("L102").Value = 3403.6667
$ws.Range("M102").Value = -455.3332999999998
$ws.Range("N102").Value = -6647.6667

$ws.Range("H136").Value = 1239.9791
$ws.Range("I136").Value = 1236.081
$ws.Range("J136").Value = 1253.091
$ws.Range("K136").Value = 3708.242999999999
$ws.Range("L136").Value = 3759.273
$ws.Range("M136").Value = -1158.242999999999
$ws.Range("N136").Value = -8859.272999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 155.96
$ws.Range("I80").Value = 69
$ws.Range("J80").Value = 196.88235
$ws.Range("K80").Value = 69
$ws.Range("L80").Value = 196.88235
$ws.Range("M80").Value = 929
$ws.Range("N80").Value = -2192.88235

$ws.Range("H83").Value = 155.96
$ws.Range("I83").Value = 69
$ws.Range("J83").Value = 196.88235
$ws.Range("K83").Value = 345
$ws.Range("L83").Value = 984.41175
$ws.Range("M83").Value = 4647
$ws.Range("N83").Value = -10968.41175

$ws.Range("H94").Value = 1348.1666
$ws.Range("I94").Value = 1394.75
$ws.Range("J94").Value = 1255
$ws.Range("K94").Value = 1394.75
$ws.Range("L94").Value = 1255
$ws.Range("M94").Value = -943.75
$ws.Range("N94").Value = -2157

$ws.Range("H99").Value = 2789.5417
$ws.Range("I99").Value = 1129.8182
$ws.Range("J99").Value = 4193.923
$ws.Range("K99").Value = 1129.8182
$ws.Range("L99").Value = 4193.923
$ws.Range("M99").Value = 368.1818000000001
$ws.Range("N99").Value = -7189.923

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 441.8
$ws.Range("I7").Value = 357
$ws.Range("J7").Value = 675
$ws.Range("K7").Value = 357
$ws.Range("L7").Value = 675
$ws.Range("M7").Value = -244
$ws.Range("N7").Value = -901

$ws.Range("H31").Value = 2159.8857
$ws.Range("I31").Value = 878.087
$ws.Range("J31").Value = 4616.6665
$ws.Range("K31").Value = 878.087
$ws.Range("L31").Value = 4616.6665
$ws.Range("M31").Value = -583.087
$ws.Range("N31").Value = -5206.6665

$ws.Range("H34").Value = 2159.8857
$ws.Range("I34").Value = 878.087
$ws.Range("J34").Value = 4616.6665
$ws.Range("K34").Value = 878.087
$ws.Range("L34").Value = 4616.6665
$ws.Range("M34").Value = -676.087
$ws.Range("N34").Value = -5020.6665

$ws.Range("H88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("N88").ClearContents()

$ws.Range("H91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("N91").ClearContents()

$ws.Range("H134").Value = 10298.385
$ws.Range("I134").Value = 25708.75
$ws.Range("J134").Value = 3449.3333
$ws.Range("K134").Value = 77126.25
$ws.Range("L134").Value = 10347.9999
$ws.Range("M134").Value = -74591.25
$ws.Range("N134").Value = -15417.9999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H69").Value = 6250

$ws.Range("H72").Value = 6250

$ws.Range("H132").Value = 3245
$ws.Range("I132").Value = 800
$ws.Range("J132").Value = 4060
$ws.Range("K132").Value = 7200
$ws.Range("L132").Value = 36540
$ws.Range("M132").Value = -4670
$ws.Range("N132").Value = -41600

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2551.6128
$ws.Range("I132").Value = 1944.1875
$ws.Range("J132").Value = 3199.5334
$ws.Range("K132").Value = 5832.5625
$ws.Range("L132").Value = 9598.600199999999
$ws.Range("M132").Value = -3302.5625
$ws.Range("N132").Value = -14658.6002

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3083.7917
$ws.Range("I7").Value = 1897.0714
$ws.Range("J7").Value = 4745.2
$ws.Range("K7").Value = 1897.0714
$ws.Range("L7").Value = 4745.2
$ws.Range("M7").Value = -1785.0714
$ws.Range("N7").Value = -4969.2

$ws.Range("H22").Value = 64581.562
$ws.Range("I22").Value = 112455.664
$ws.Range("J22").Value = 3029.1428
$ws.Range("K22").Value = 112455.664
$ws.Range("L22").Value = 3029.1428
$ws.Range("M22").Value = -112160.664
$ws.Range("N22").Value = -3619.1428

$ws.Range("H27").Value = 64581.562
$ws.Range("I27").Value = 112455.664
$ws.Range("J27").Value = 3029.1428
$ws.Range("K27").Value = 112455.664
$ws.Range("L27").Value = 3029.1428
$ws.Range("M27").Value = -112348.664
$ws.Range("N27").Value = -3243.1428

$ws.Range("H40").Value = 6574.4614
$ws.Range("I40").Value = 5956.3
$ws.Range("J40").Value = 8635
$ws.Range("K40").Value = 5956.3
$ws.Range("L40").Value = 8635
$ws.Range("M40").Value = -5820.3
$ws.Range("N40").Value = -8907

$ws.Range("H42").Value = 34793.4
$ws.Range("J42").Value = 34793.4
$ws.Range("L42").Value = 34793.4
$ws.Range("N42").Value = -35919.4

$ws.Range("H49").Value = 34793.4
$ws.Range("J49").Value = 34793.4
$ws.Range("L49").Value = 34793.4
$ws.Range("N49").Value = -35087.4

$ws.Range("H64").Value = 33760
$ws.Range("J64").Value = 33760
$ws.Range("L64").Value = 33760
$ws.Range("N64").Value = -34210

$ws.Range("H67").Value = 33760
$ws.Range("J67").Value = 33760
$ws.Range("L67").Value = 33760
$ws.Range("N67").Value = -35320

$ws.Range("H122").Value = 3009.3333
$ws.Range("I122").Value = 1711.2
$ws.Range("J122").Value = 9500
$ws.Range("K122").Value = 5133.6
$ws.Range("L122").Value = 28500
$ws.Range("M122").Value = -2683.6
$ws.Range("N122").Value = -33400

$ws.Range("H126").Value = 3083.7917
$ws.Range("I126").Value = 1897.0714
$ws.Range("J126").Value = 4745.2
$ws.Range("K126").Value = 5691.2142
$ws.Range("L126").Value = 14235.6
$ws.Range("M126").Value = -3221.2142
$ws.Range("N126").Value = -19175.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 27780944
$ws.Range("I132").Value = 1622.25
$ws.Range("J132").Value = 83339580
$ws.Range("K132").Value = 4866.75
$ws.Range("L132").Value = 250018740
$ws.Range("M132").Value = -2336.75
$ws.Range("N132").Value = -250023800
